# Backlog update: add RF31..RF34 (commit: "Atualização do backlog (adicionados RF31 ao RF34)")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend the table with 4 new rows (32-35), reusing the existing
# --- alternating row styles (row 30 -> style "Importante" fill, row 31 ->
# --- style "Desejavel" fill) so the new rows keep the same banding.
$ws.Range("A30:D30").Copy()
$ws.Range("A32:D32").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A34:D34").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("A31:D31").Copy()
$ws.Range("A33:D33").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A35:D35").PasteSpecial(-4122)   # xlPasteFormats

$ws.Rows.Item(32).RowHeight = 20.25
$ws.Rows.Item(33).RowHeight = 20.25
$ws.Rows.Item(34).RowHeight = 20.25
$ws.Rows.Item(35).RowHeight = 20.25

# Column A (Requisito id) filled first for every new row ...
$ws.Range("A32").Value = "RF31"
$ws.Range("A33").Value = "RF32"
$ws.Range("A34").Value = "RF33"
$ws.Range("A35").Value = "RF34"

# ... then Requisito / Descricao filled row by row ...
$ws.Range("B32").Value = "Diagrama de Visão de Negócio"
$ws.Range("C32").Value = "Criar e adicionar como parte da documentação"

$ws.Range("B33").Value = "Sequencia de Fibonacci"
$ws.Range("C33").Value = "Adicionar ao backlog, para criar o gráfico de Burndown"

$ws.Range("B34").Value = "Gráfico de Burndown"
$ws.Range("C34").Value = "Incluir ao backlog do produto"

$ws.Range("B35").Value = "Service Level Agreement do projeto"
$ws.Range("C35").Value = "Criar o SLA -  responsabilidade entre cliente e provedor"

# ... and finally the Classificação column.
$ws.Range("D32").Value = "Importante"
$ws.Range("D33").Value = "Importante"
$ws.Range("D34").Value = "Importante"
$ws.Range("D35").Value = "Desejável"

# Match the saved view state: zoomed out a bit, and the cursor resting on
# the first empty row right below the table that was just extended.
$win = $wb.Windows.Item(1)
$win.Zoom = 55
$ws.Range("D36").Select() | Out-Null
